$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 updates
$ws.Range("A7").Value2 = 131106436
$ws.Range("B7").Value2 = 5493
$ws.Range("E7").Value2 = 101410
$ws.Range("F7").Value2 = "Reliktbock"
$ws.Range("G7").Value2 = "Nothorhina muricata"
$ws.Range("H7").Value2 = "(Dalman, 1817)"
$ws.Range("I7").Value2 = "2"
$ws.Range("J7").Value2 = "ex."
$ws.Range("P7").Value2 = "Svartmyran, Mpd"
$ws.Range("Q7").Value2 = 616762
$ws.Range("R7").Value2 = 6934714
$ws.Range("X7").Value2 = "2025_0743"
$ws.Range("Z7").Value2 = "11:39"
$ws.Range("AB7").Value2 = "11:39"
$ws.Range("AC7").Value2 = "Två kläckhål"
$ws.Range("AX7").Value2 = "David Isaksson"

# Row 8 updates
$ws.Range("A8").Value2 = 131108352
$ws.Range("B8").Value2 = 80216
$ws.Range("E8").Value2 = 388
$ws.Range("F8").Value2 = "Stiftgelélav"
$ws.Range("G8").Value2 = "Collema furfuraceum"
$ws.Range("H8").Value2 = "(Arnold) Du Rietz"
$ws.Range("I8").Value2 = "1"
$ws.Range("J8").Value2 = "bålar"
$ws.Range("P8").Value2 = "S Svartmyran, Mpd"
$ws.Range("Q8").Value2 = 616863
$ws.Range("R8").Value2 = 6934788
$ws.Range("X8").Value2 = "2025_0758"
$ws.Range("Z8").Value2 = "14:47"
$ws.Range("AB8").Value2 = "14:47"
$ws.Range("AX8").Value2 = "Måns Svensson"
